$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 text update (company name) ---
$ws.Range("B3").Value = "Egypt Gas Company SAE (CASE:EGAS)"

# --- Column D update for both rows ---
$ws.Range("D2").Value = 0.261
$ws.Range("D3").Value = 0.261

# --- Column E is removed entirely for both rows ---
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()

# --- Remaining numeric columns updated identically for rows 2 and 3 ---
$updates = @{
    "G"  = -0.05068143100511073
    "H"  = -0.05068143100511073
    "I"  = -0.0706984667802385
    "J"  = -0.0565587734241908
    "K"  = 9.27
    "L"  = 0.03948040885860307
    "M"  = 1.41
    "N"  = 0.01060948081264108
    "O"  = 0.1521035598705502
    "P"  = 1.41
    "Q"  = 0.01060948081264108
    "R"  = 0.1521035598705502
    "U"  = 11.9
    "V"  = 0.08954100827689992
    "W"  = 0.1459842519685039
    "X"  = 0.08581654877709093
    "Y"  = 0.060167703191413
    "Z"  = 2.078060005310205
    "AA" = -0.1175325250022126
    "AB" = 0.07289417923711761
    "AC" = -0.1904267042393302
    "AD" = 84.5
    "AF" = 84.5
    "AG" = 72.59999999999999
    "AH" = 0.3886844526218951
    "AI" = 0.3787539220080681
    "AJ" = 0.3532846715328467
    "AK" = 0.34375
    "AM" = -27
    "AN" = -6.550387596899225
    "AP" = -5.627906976744185
    "AQ" = 0.6148148148148148
}

foreach ($col in $updates.Keys) {
    $value = $updates[$col]
    $ws.Range("$col`2").Value = $value
    $ws.Range("$col`3").Value = $value
}
